$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets("ALC")
$ws.Range("H51").Value = 7080.7617
$ws.Range("I51").Value = 10853.546
$ws.Range("J51").Value = 2930.7
$ws.Range("K51").Value = 10853.546
$ws.Range("L51").Value = 2930.7
$ws.Range("M51").Value = -10369.546
$ws.Range("N51").Value = -3898.7
$ws.Range("H87").Value = 31164.6
$ws.Range("J87").Value = 31164.6
$ws.Range("L87").Value = 31164.6
$ws.Range("N87").Value = -33660.6
$ws.Range("H90").Value = 31164.6
$ws.Range("J90").Value = 31164.6
$ws.Range("L90").Value = 93493.79999999999
$ws.Range("N90").Value = -105973.8
$ws.Range("H137").Value = 1390.0212
$ws.Range("I137").Value = 1368.4615
$ws.Range("J137").Value = 1398.2646
$ws.Range("K137").Value = 4105.3845
$ws.Range("L137").Value = 4194.793799999999
$ws.Range("M137").Value = -1555.3845
$ws.Range("N137").Value = -9294.793799999999
$ws.Range("H138").Value = 4293.8438
$ws.Range("I138").Value = 2248.2
$ws.Range("J138").Value = 5223.6816
$ws.Range("K138").Value = 6744.599999999999
$ws.Range("L138").Value = 15671.0448
$ws.Range("M138").Value = -1604.599999999999
$ws.Range("N138").Value = -25951.0448

# ---- ARM ----
$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 32353.697
$ws.Range("I32").Value = 5539.533
$ws.Range("J32").Value = 142048
$ws.Range("K32").Value = 5539.533
$ws.Range("L32").Value = 142048
$ws.Range("M32").Value = -5252.533
$ws.Range("N32").Value = -142622
$ws.Range("H74").Value = 1634.2826
$ws.Range("I74").Value = 953.0417
$ws.Range("J74").Value = 2377.4546
$ws.Range("K74").Value = 953.0417
$ws.Range("L74").Value = 2377.4546
$ws.Range("M74").Value = -79.04169999999999
$ws.Range("N74").Value = -4125.4546
$ws.Range("H75").Value = 5000
$ws.Range("I75").Value = 5000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 5000
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0
$ws.Range("M75").Value = -4126
$ws.Range("H77").Value = 1634.2826
$ws.Range("I77").Value = 953.0417
$ws.Range("J77").Value = 2377.4546
$ws.Range("K77").Value = 4765.2085
$ws.Range("L77").Value = 11887.273
$ws.Range("M77").Value = -397.2084999999997
$ws.Range("N77").Value = -20623.273
$ws.Range("H78").Value = 5000
$ws.Range("I78").Value = 5000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 15000
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0
$ws.Range("M78").Value = -10632

# ---- BSM ----
$ws = $wb.Worksheets("BSM")
$ws.Range("H99").Value = 1663.7646
$ws.Range("J99").Value = 1653.3846
$ws.Range("L99").Value = 1653.3846
$ws.Range("N99").Value = -4649.3846
$ws.Range("H118").Value = 26300
$ws.Range("J118").Value = 26300
$ws.Range("L118").Value = 26300
$ws.Range("N118").Value = -29614

# ---- CRP ----
$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 17574.484
$ws.Range("I31").Value = 46540.5
$ws.Range("J31").Value = 2401.8096
$ws.Range("K31").Value = 46540.5
$ws.Range("L31").Value = 2401.8096
$ws.Range("M31").Value = -46245.5
$ws.Range("N31").Value = -2991.8096
$ws.Range("H34").Value = 17574.484
$ws.Range("I34").Value = 46540.5
$ws.Range("J34").Value = 2401.8096
$ws.Range("K34").Value = 46540.5
$ws.Range("L34").Value = 2401.8096
$ws.Range("M34").Value = -46338.5
$ws.Range("N34").Value = -2805.8096
$ws.Range("H68").Value = 13162
$ws.Range("J68").Value = 13162
$ws.Range("L68").Value = 13162
$ws.Range("N68").Value = -14660
$ws.Range("H71").Value = 13162
$ws.Range("J71").Value = 13162
$ws.Range("L71").Value = 39486
$ws.Range("N71").Value = -46974
$ws.Range("H74").Value = 21480.25
$ws.Range("J74").Value = 21480.25
$ws.Range("L74").Value = 21480.25
$ws.Range("N74").Value = -23228.25
$ws.Range("H77").Value = 21480.25
$ws.Range("J77").Value = 21480.25
$ws.Range("L77").Value = 64440.75
$ws.Range("N77").Value = -73176.75
$ws.Range("H108").Value = 29800
$ws.Range("J108").Value = 29800
$ws.Range("L108").Value = 29800
$ws.Range("N108").Value = -37480
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0

# ---- CUL ----
$ws = $wb.Worksheets("CUL")
$ws.Range("H37").Value = 1449718.1
$ws.Range("J37").Value = 1449718.1
$ws.Range("L37").Value = 4349154.300000001
$ws.Range("N37").Value = -4349378.300000001
$ws.Range("H39").Value = 3433.3333
$ws.Range("I39").Value = 3400
$ws.Range("J39").Value = 3450
$ws.Range("K39").Value = 10200
$ws.Range("L39").Value = 10350
$ws.Range("M39").Value = -9906
$ws.Range("N39").Value = -10938
$ws.Range("H113").Value = 544.3714
$ws.Range("I113").Value = 537.6923
$ws.Range("J113").Value = 548.3182
$ws.Range("K113").Value = 1613.0769
$ws.Range("L113").Value = 1644.9546
$ws.Range("M113").Value = 556.9231
$ws.Range("N113").Value = -5984.9546
$ws.Range("H122").Value = 7998.0713
$ws.Range("I122").Value = 261.33334
$ws.Range("J122").Value = 10108.091
$ws.Range("K122").Value = 2352.00006
$ws.Range("L122").Value = 90972.819
$ws.Range("M122").Value = 97.9999399999997
$ws.Range("N122").Value = -95872.819
$ws.Range("H123").Value = 3403.625
$ws.Range("I123").Value = 1882.25
$ws.Range("J123").Value = 4925
$ws.Range("K123").Value = 5646.75
$ws.Range("L123").Value = 14775
$ws.Range("M123").Value = -3196.75
$ws.Range("N123").Value = -19675
$ws.Range("H129").Value = 1875.5333
$ws.Range("I129").Value = 533.3333
$ws.Range("J129").Value = 2770.3333
$ws.Range("K129").Value = 1599.9999
$ws.Range("L129").Value = 8310.999899999999
$ws.Range("M129").Value = 3400.0001
$ws.Range("N129").Value = -18310.9999
$ws.Range("H131").Value = 1254.0146
$ws.Range("I131").Value = 865
$ws.Range("J131").Value = 1321.0862
$ws.Range("K131").Value = 2595
$ws.Range("L131").Value = 3963.2586
$ws.Range("M131").Value = 2445
$ws.Range("N131").Value = -14043.2586

# ---- GSM ----
$ws = $wb.Worksheets("GSM")
$ws.Range("H64").Value = 41262
$ws.Range("J64").Value = 41262
$ws.Range("L64").Value = 41262
$ws.Range("N64").Value = -41758
$ws.Range("H67").Value = 41262
$ws.Range("J67").Value = 41262
$ws.Range("L67").Value = 41262
$ws.Range("N67").Value = -42978
$ws.Range("H102").Value = 302017.7
$ws.Range("I102").Value = 2130.8333
$ws.Range("K102").Value = 2130.8333
$ws.Range("M102").Value = -508.8332999999998

# ---- LTW ----
$ws = $wb.Worksheets("LTW")
$ws.Range("H40").Value = 501250
$ws.Range("I40").Value = 1000000
$ws.Range("J40").Value = 2500
$ws.Range("K40").Value = 1000000
$ws.Range("L40").Value = 2500
$ws.Range("M40").Value = -999864
$ws.Range("N40").Value = -2772
$ws.Range("H68").Value = 2335.6191
$ws.Range("I68").Value = 1369.6666
$ws.Range("J68").Value = 4750.5
$ws.Range("K68").Value = 1369.6666
$ws.Range("L68").Value = 4750.5
$ws.Range("M68").Value = -620.6666
$ws.Range("N68").Value = -6248.5
$ws.Range("H71").Value = 2335.6191
$ws.Range("I71").Value = 1369.6666
$ws.Range("J71").Value = 4750.5
$ws.Range("K71").Value = 6848.333000000001
$ws.Range("L71").Value = 23752.5
$ws.Range("M71").Value = -3104.333000000001
$ws.Range("N71").Value = -31240.5
$ws.Range("H132").Value = 2779.121
$ws.Range("I132").Value = 2954.0833
$ws.Range("J132").Value = 2312.5557
$ws.Range("K132").Value = 8862.249899999999
$ws.Range("L132").Value = 6937.6671
$ws.Range("M132").Value = -6332.249899999999
$ws.Range("N132").Value = -11997.6671

# ---- WVR ----
$ws = $wb.Worksheets("WVR")
$ws.Range("H21").Value = 12000
$ws.Range("J21").Value = 12000
$ws.Range("L21").Value = 12000
$ws.Range("N21").Value = -12470
$ws.Range("H35").Value = 12000
$ws.Range("J35").Value = 12000
$ws.Range("L35").Value = 12000
$ws.Range("N35").Value = -12580
$ws.Range("H126").Value = 1243.6
$ws.Range("I126").Value = 1341.2222
$ws.Range("J126").Value = 1097.1666
$ws.Range("K126").Value = 4023.6666
$ws.Range("L126").Value = 3291.4998
$ws.Range("M126").Value = -1553.6666
$ws.Range("N126").Value = -8231.4998
$ws.Range("H136").Value = 940.4
$ws.Range("I136").Value = 500.57895
$ws.Range("J136").Value = 2333.1667
$ws.Range("K136").Value = 1501.73685
$ws.Range("L136").Value = 6999.500100000001
$ws.Range("M136").Value = 1048.26315
$ws.Range("N136").Value = -12099.5001
